# Sort the curvature calibration data (rows 2-8, columns A-D) into
# ascending chronological order based on column A (time).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nRows = 7
$nCols = 4
$firstRow = 2

# Read current data rows 2-8 (A:D) cell-by-cell into an array of arrays.
$rows = @()
for ($i = 0; $i -lt $nRows; $i++) {
    $r = $firstRow + $i
    $row = @(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2
    )
    $rows += ,$row
}

# Sort rows ascending by column A (index 0).
$sortedRows = $rows | Sort-Object { $_[0] }

# Write the sorted rows back into A2:D8.
for ($i = 0; $i -lt $nRows; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 1).Value2 = $sortedRows[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $sortedRows[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $sortedRows[$i][2]
    $ws.Cells.Item($r, 4).Value2 = $sortedRows[$i][3]
}
